$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update attendance marks for the removed date (column U) ---
# Column U (one date's P/A column) is being retired: clear its marks for
# every participant row (7-82) and let the adjacent blank column (V)
# lend its format, matching how the rest of the unused date columns
# already look.
$ws.Range("V7:V82").Copy()
$ws.Range("U7:U82").PasteSpecial(-4122)
$ws.Range("U7:U82").ClearContents()

# Row 61 (S.No. 55) had their column T attendance mark corrected from
# Present to Absent.
$ws.Range("T61").Value = "A"

# --- Re-arrange the freeze panes: freeze only the header rows (1-6),
# no longer also freezing the first three columns. ---
$ws.Range("A7").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B8").Select()
